$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column Q is entirely empty in this sheet; delete it so that the
# data that was in columns R:AH shifts left to occupy Q:AG.
$ws.Columns("Q").Delete()
